$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.068.79"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "'1.816.58"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").Value = "'233.00"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'0.5866"
$ws.Range("E6").Value = "  -3.29%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").Value = "'0.2714"
$ws.Range("E8").Value = "  -4.01%  "
$ws.Range("D9").Value = "'0.06757"
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("D10").Value = "'22.76"
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "'1.815.65"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "'4.625"
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").Value = "'0.6166"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "'0.000009388"
$ws.Range("E15").Value = "  -5.84%  "
$ws.Range("D16").Value = "'74.50"
$ws.Range("E16").Value = "  -6.65%  "
$ws.Range("D17").Value = "'28.864.82"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "'5.403"
$ws.Range("E18").Value = "  -9.91%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("D20").Value = "'206.83"
$ws.Range("E20").Value = "  -10.44%  "
$ws.Range("D21").Value = "'11.38"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").Value = "'6.730"
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "'154.01"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "'7.757"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").Value = "'0.1253"
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").Value = "'16.18"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("D28").Value = "'1.406"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "'0.06274"
$ws.Range("E29").Value = "  -5.43%  "
$ws.Range("D30").Value = "'1.428"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").Value = "'3.688"
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").Value = "'3.679"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").Value = "'1.681"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("D34").Value = "'1.045"
$ws.Range("E34").Value = "  -7.64%  "
$ws.Range("D35").Value = "'2.533"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").Value = "'0.6323"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("D37").Value = "'2.745"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "'0.01701"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").Value = "'6.373"
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "'1.128.76"
$ws.Range("E40").Value = "  -8.57%  "
$ws.Range("D41").Value = "'0.8579"
$ws.Range("E41").Value = "  -7.52%  "
$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("B43").Value = "'Quant"
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.92"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("B44").Value = "'RocketPoolETH"
$ws.Range("C44").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "'1.962.07"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'59.98"
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "'0.4539"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'0.05484"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "'1.554"
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("D50").Value = "'8.199"
$ws.Range("E50").Value = "  -3.86%  "
$ws.Range("E51").Value = "  -0.23%  "
